$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -7
$ws.Range("F4").Value = 7
$ws.Range("F9").Value = -13
$ws.Range("F11").Value = 4
$ws.Range("F12").Value = -4
$ws.Range("F13").Value = -3
$ws.Range("F15").Value = -1
$ws.Range("F19").Value = 2
